$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 1708.1666
$ws.Range("J6").Value = 0
$ws.Range("L6").Value = 0
$ws.Range("N6").ClearContents()
$ws.Range("H125").Value = 1833.3334
$ws.Range("I125").Value = 1000
$ws.Range("K125").Value = 9000
$ws.Range("M125").Value = -6540
$ws.Range("H132").Value = 1297.0714
$ws.Range("I132").Value = 1297.0714
$ws.Range("K132").Value = 3891.2142
$ws.Range("M132").Value = -1361.2142
$ws.Range("H137").Value = 1672.7222
$ws.Range("I137").Value = 1449.0465
$ws.Range("J137").Value = 2547.0908
$ws.Range("K137").Value = 4347.139499999999
$ws.Range("L137").Value = 7641.2724
$ws.Range("M137").Value = -1797.139499999999
$ws.Range("N137").Value = -12741.2724
$ws.Range("H138").Value = 4604.643
$ws.Range("I138").Value = 4053.0908
$ws.Range("J138").Value = 4961.5293
$ws.Range("K138").Value = 12159.2724
$ws.Range("L138").Value = 14884.5879
$ws.Range("M138").Value = -7019.2724
$ws.Range("N138").Value = -25164.5879

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 1733.3334
$ws.Range("I4").Value = 100
$ws.Range("K4").Value = 100
$ws.Range("M4").Value = 16
$ws.Range("H44").Value = 39995
$ws.Range("J44").Value = 39995
$ws.Range("L44").Value = 39995
$ws.Range("N44").Value = -40971
$ws.Range("H45").Value = 3111.1538
$ws.Range("I45").Value = 2725.7778
$ws.Range("K45").Value = 2725.7778
$ws.Range("M45").Value = -2348.7778
$ws.Range("H61").Value = 1731.7097
$ws.Range("I61").Value = 1339.15
$ws.Range("K61").Value = 1339.15
$ws.Range("M61").Value = -1127.15
$ws.Range("H74").Value = 1445.6052
$ws.Range("I74").Value = 1001.4667
$ws.Range("J74").Value = 3111.125
$ws.Range("K74").Value = 1001.4667
$ws.Range("L74").Value = 3111.125
$ws.Range("M74").Value = -127.4666999999999
$ws.Range("N74").Value = -4859.125
$ws.Range("H77").Value = 1445.6052
$ws.Range("I77").Value = 1001.4667
$ws.Range("J77").Value = 3111.125
$ws.Range("K77").Value = 5007.3335
$ws.Range("L77").Value = 15555.625
$ws.Range("M77").Value = -639.3334999999997
$ws.Range("N77").Value = -24291.625
$ws.Range("H122").Value = 4199
$ws.Range("I122").Value = 3998.75
$ws.Range("J122").Value = 5000
$ws.Range("K122").Value = 11996.25
$ws.Range("L122").Value = 15000
$ws.Range("M122").Value = -9546.25
$ws.Range("N122").Value = -19900
$ws.Range("H132").Value = 2650.2942
$ws.Range("I132").Value = 1588.25
$ws.Range("K132").Value = 4764.75
$ws.Range("M132").Value = -2234.75
$ws.Range("H136").Value = 1731.7097
$ws.Range("I136").Value = 1339.15
$ws.Range("K136").Value = 4017.45
$ws.Range("M136").Value = -1467.45

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 843.1667
$ws.Range("I94").Value = 744.4
$ws.Range("K94").Value = 744.4
$ws.Range("M94").Value = -293.4
$ws.Range("H134").Value = 2968.04
$ws.Range("I134").Value = 2800.3684
$ws.Range("J134").Value = 3499
$ws.Range("K134").Value = 8401.1052
$ws.Range("L134").Value = 10497
$ws.Range("M134").Value = -5866.1052
$ws.Range("N134").Value = -15567

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 219.42857
$ws.Range("I7").Value = 89.5
$ws.Range("K7").Value = 89.5
$ws.Range("M7").Value = 23.5
$ws.Range("H22").Value = 650
$ws.Range("J22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("N22").ClearContents()
$ws.Range("H31").Value = 2154.4348
$ws.Range("I31").Value = 2179.25
$ws.Range("K31").Value = 2179.25
$ws.Range("M31").Value = -1884.25
$ws.Range("H34").Value = 2154.4348
$ws.Range("I34").Value = 2179.25
$ws.Range("K34").Value = 2179.25
$ws.Range("M34").Value = -1977.25
$ws.Range("H42").Value = 4000
$ws.Range("J42").Value = 4000
$ws.Range("L42").Value = 4000
$ws.Range("N42").Value = -5186
$ws.Range("H122").Value = 2612.5
$ws.Range("J122").Value = 2612.5
$ws.Range("L122").Value = 7837.5
$ws.Range("N122").Value = -12737.5
$ws.Range("H132").Value = 4310.7896
$ws.Range("I132").Value = 4229.8237
$ws.Range("K132").Value = 12689.4711
$ws.Range("M132").Value = -10159.4711
$ws.Range("H141").Value = 61555
$ws.Range("J141").Value = 61555
$ws.Range("L141").Value = 61555
$ws.Range("N141").Value = -71915

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H120").Value = 5465
$ws.Range("I120").Value = 5465
$ws.Range("K120").Value = 16395
$ws.Range("M120").Value = -11557
$ws.Range("H122").Value = 168749.5
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 168749.5
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 1518745.5
$ws.Range("M122").ClearContents()
$ws.Range("N122").Value = -1523645.5
$ws.Range("H129").Value = 1471
$ws.Range("J129").Value = 1782.8334
$ws.Range("L129").Value = 5348.5002
$ws.Range("N129").Value = -15348.5002
$ws.Range("H131").Value = 1951
$ws.Range("J131").Value = 1997.25
$ws.Range("L131").Value = 5991.75
$ws.Range("N131").Value = -16071.75
$ws.Range("H132").Value = 3781.8333
$ws.Range("I132").Value = 3599.8
$ws.Range("J132").Value = 3851.8462
$ws.Range("K132").Value = 32398.2
$ws.Range("L132").Value = 34666.6158
$ws.Range("M132").Value = -29868.2
$ws.Range("N132").Value = -39726.6158
$ws.Range("H138").Value = 5639.4
$ws.Range("I138").Value = 5639.4
$ws.Range("K138").Value = 16918.2
$ws.Range("M138").Value = -11778.2

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 82.5
$ws.Range("I2").Value = 90.71429000000001
$ws.Range("K2").Value = 90.71429000000001
$ws.Range("M2").Value = 22.28570999999999
$ws.Range("H132").Value = 4550.684
$ws.Range("I132").Value = 4404.857
$ws.Range("J132").Value = 4959
$ws.Range("K132").Value = 13214.571
$ws.Range("L132").Value = 14877
$ws.Range("M132").Value = -10684.571
$ws.Range("N132").Value = -19937

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 0
$ws.Range("I40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("M40").ClearContents()
$ws.Range("H43").Value = 428565
$ws.Range("J43").Value = 428565
$ws.Range("L43").Value = 428565
$ws.Range("N43").Value = -428951
$ws.Range("H93").Value = 2058
$ws.Range("I93").Value = 2002.8462
$ws.Range("K93").Value = 2002.8462
$ws.Range("M93").Value = -754.8462
$ws.Range("H122").Value = 5925
$ws.Range("I122").Value = 850
$ws.Range("J122").Value = 11000
$ws.Range("K122").Value = 2550
$ws.Range("L122").Value = 33000
$ws.Range("M122").Value = -100
$ws.Range("N122").Value = -37900
$ws.Range("H132").Value = 5750.364
$ws.Range("I132").Value = 5636.7144
$ws.Range("K132").Value = 16910.1432
$ws.Range("M132").Value = -14380.1432
$ws.Range("H136").Value = 5714.722
$ws.Range("I136").Value = 5791.2
$ws.Range("J136").Value = 5332.3335
$ws.Range("K136").Value = 17373.6
$ws.Range("L136").Value = 15997.0005
$ws.Range("M136").Value = -14823.6
$ws.Range("N136").Value = -21097.0005

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 63881.75
$ws.Range("I62").Value = 84559.336
$ws.Range("K62").Value = 84559.336
$ws.Range("M62").Value = -83935.336
$ws.Range("H65").Value = 63881.75
$ws.Range("I65").Value = 84559.336
$ws.Range("K65").Value = 422796.68
$ws.Range("M65").Value = -419676.68
$ws.Range("H103").Value = 25069
$ws.Range("I103").Value = 25069
$ws.Range("K103").Value = 25069
$ws.Range("M103").Value = -23897
$ws.Range("H113").Value = 1510.4546
$ws.Range("I113").Value = 2137.8572
$ws.Range("J113").Value = 412.5
$ws.Range("K113").Value = 6413.571599999999
$ws.Range("L113").Value = 1237.5
$ws.Range("M113").Value = -4243.571599999999
$ws.Range("N113").Value = -5577.5
$ws.Range("H132").Value = 1761.5238
$ws.Range("I132").Value = 1287.5333
$ws.Range("J132").Value = 2946.5
$ws.Range("K132").Value = 3862.5999
$ws.Range("L132").Value = 8839.5
$ws.Range("M132").Value = -1332.5999
$ws.Range("N132").Value = -13899.5
$ws.Range("H136").Value = 1670.4166
$ws.Range("I136").Value = 1289.125
$ws.Range("J136").Value = 2433
$ws.Range("K136").Value = 3867.375
$ws.Range("L136").Value = 7299
$ws.Range("M136").Value = -1317.375
$ws.Range("N136").Value = -12399
